$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted at row 17 (pushing the
# existing rows 17-88 down to 18-89, which is exactly what the diff shows:
# every row's data equals the row above it from the previous version).
$ws.Rows.Item(17).Insert()

$ws.Range("A17").Value = 8
$ws.Range("B17").Value = "Terminal La Palmera de La Serena"
$ws.Range("C17").Value = "Coquimbo"
$ws.Range("D17").Value = 44701
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 100112052
$ws.Range("G17").Value = "Albahaca"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 1120
$ws.Range("K17").Value = 4000
$ws.Range("L17").Value = 4500
$ws.Range("M17").Value = 4250
$ws.Range("N17").Value = "$/paquete"
$ws.Range("O17").Value = "Región de Arica y Parinacota"
$ws.Range("P17").Value = 4250
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = "Hortaliza"
